$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2468
$ws.Range("I43").Value = 1365.6666
$ws.Range("J43").Value = 2881.375
$ws.Range("K43").Value = 1365.6666
$ws.Range("L43").Value = 2881.375
$ws.Range("M43").Value = -1296.6666
$ws.Range("N43").Value = -3019.375
$ws.Range("H53").Value = 1000.6667
$ws.Range("I53").Value = 1000.6667
$ws.Range("K53").Value = 1000.6667
$ws.Range("M53").Value = -363.6667
$ws.Range("H98").Value = 1519.7693
$ws.Range("I98").Value = 978.2727
$ws.Range("J98").Value = 4498
$ws.Range("K98").Value = 978.2727
$ws.Range("L98").Value = 4498
$ws.Range("M98").Value = 519.7273
$ws.Range("N98").Value = -7494
$ws.Range("H101").Value = 693.3333
$ws.Range("I101").Value = 889.125
$ws.Range("K101").Value = 2667.375
$ws.Range("M101").Value = -1045.375
$ws.Range("H120").Value = 69846
$ws.Range("J120").Value = 69846
$ws.Range("L120").Value = 69846
$ws.Range("N120").Value = -79522
$ws.Range("H122").Value = 1519.7693
$ws.Range("I122").Value = 978.2727
$ws.Range("J122").Value = 4498
$ws.Range("K122").Value = 2934.8181
$ws.Range("L122").Value = 13494
$ws.Range("M122").Value = -484.8181
$ws.Range("N122").Value = -18394
$ws.Range("H138").Value = 3739.039
$ws.Range("I138").Value = 3623.0286
$ws.Range("J138").Value = 3835.7144
$ws.Range("K138").Value = 10869.0858
$ws.Range("L138").Value = 11507.1432
$ws.Range("M138").Value = -5729.085800000001
$ws.Range("N138").Value = -21787.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 25029250
$ws.Range("I34").Value = 100000000
$ws.Range("K34").Value = 100000000
$ws.Range("M34").Value = -99999729

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2333.3333
$ws.Range("I94").Value = 779.4783
$ws.Range("K94").Value = 779.4783
$ws.Range("M94").Value = -328.4783

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 19999
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2923.5
$ws.Range("I34").Value = 1950.5
$ws.Range("J34").Value = 3247.8333
$ws.Range("K34").Value = 5851.5
$ws.Range("L34").Value = 9743.499899999999
$ws.Range("M34").Value = -5767.5
$ws.Range("N34").Value = -9911.499899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2032.5
$ws.Range("I132").Value = 2035.5454
$ws.Range("K132").Value = 6106.6362
$ws.Range("M132").Value = -3576.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 2531500
$ws.Range("J43").Value = 5024250
$ws.Range("L43").Value = 5024250
$ws.Range("N43").Value = -5024636
$ws.Range("H124").Value = 67500
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 67500
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 67500
$ws.Range("N124").Value = -77320
$ws.Range("H125").Value = 67499.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 67499.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 67499.5
$ws.Range("N125").Value = -77339.5
$ws.Range("H127").Value = 67500
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 67500
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 67500
$ws.Range("N127").Value = -77420
$ws.Range("H128").Value = 71374.75
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 71374.75
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 71374.75
$ws.Range("N128").Value = -81334.75
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H130").Value = 91500
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 91500
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 91500
$ws.Range("N130").Value = -101540
$ws.Range("H131").Value = 89500
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 89500
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 89500
$ws.Range("N131").Value = -99580
$ws.Range("H132").Value = 2799.647
$ws.Range("I132").Value = 2772.9333
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8318.7999
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -5788.7999
$ws.Range("N132").Value = -14060
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 81643
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 81643
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 81643
$ws.Range("N135").Value = -91783
$ws.Range("H136").Value = 2006.7693
$ws.Range("I136").Value = 1711.7826
$ws.Range("J136").Value = 4268.3335
$ws.Range("K136").Value = 5135.3478
$ws.Range("L136").Value = 12805.0005
$ws.Range("M136").Value = -2585.3478
$ws.Range("N136").Value = -17905.0005
$ws.Range("H137").Value = 110000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 110000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 110000
$ws.Range("N137").Value = -120200
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 67499
$ws.Range("I139").Value = 67499
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 67499
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -62359
$ws.Range("H140").Value = 65000
$ws.Range("I140").Value = 65000
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 65000
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -59820
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H62").Value = 7166.3335
$ws.Range("I62").Value = 7249.5
$ws.Range("K62").Value = 7249.5
$ws.Range("M62").Value = -6625.5
$ws.Range("H65").Value = 7166.3335
$ws.Range("I65").Value = 7249.5
$ws.Range("K65").Value = 36247.5
$ws.Range("M65").Value = -33127.5
